# Update column F (dSF) values on the active sheet with newly re-pulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = -2
    3  = -1
    4  = 1
    5  = -1
    6  = -5
    7  = -2
    8  = -1
    9  = 1
    10 = -4
    11 = 0
    12 = 0
    13 = 2
    14 = 2
    15 = 3
    16 = 0
    17 = 4
    18 = 2
    19 = -4
    20 = -4
    21 = -2
    22 = -3
    23 = 3
    24 = -2
    25 = 1
    26 = 6
    27 = 2
    28 = -3
    29 = 0
    30 = 5
    31 = -3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("F$row").Value = $newValues[$row]
}
